$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "y"
$ws.Range("C4").Value = "y"
$ws.Range("C55").Value = "y"
$ws.Range("C58").Value = "y"
$ws.Range("C100").Value = "y"
$ws.Range("D101").Value = "Query"
$ws.Range("C108").Value = "y"
$ws.Range("D117").Value = "Query"
$ws.Range("D124").Value = "Query"
$ws.Range("C136").Value = "ND"
$ws.Range("C137").Value = "ND"
$ws.Range("C149").Value = "y"
$ws.Range("C150").Value = "y"
$ws.Range("C151").Value = "y"
$ws.Range("C152").Value = "y"
$ws.Range("C153").Value = "y"
$ws.Range("C154").Value = "y"
$ws.Range("C155").Value = "y"
$ws.Range("C156").Value = "y"
$ws.Range("C157").Value = "y"
$ws.Range("C158").Value = "y"
$ws.Range("C159").Value = "y"
$ws.Range("C161").Value = "y"
$ws.Range("C162").Value = "y"
$ws.Range("C163").Value = "y"
$ws.Range("C164").Value = "y"
$ws.Range("C165").Value = "y"
$ws.Range("C166").Value = "y"
$ws.Range("D167").Value = "query"
$ws.Range("C168").Value = "y"
$ws.Range("C169").Value = "y"
$ws.Range("C170").Value = "y"
$ws.Range("C171").Value = "y"
$ws.Range("C173").Value = "y"
$ws.Range("C174").Value = "y"
$ws.Range("C175").Value = "y"
$ws.Range("C176").Value = "y"
$ws.Range("C177").Value = "y"
$ws.Range("C178").Value = "y"
$ws.Range("C179").Value = "y"
$ws.Range("C180").Value = "y"
$ws.Range("C181").Value = "y"
$ws.Range("C182").Value = "y"
$ws.Range("C183").Value = "y"
$ws.Range("C184").Value = "ND"
$ws.Range("C185").Value = "ND"
$ws.Range("C187").Value = "y"
$ws.Range("C188").Value = "y"
$ws.Range("C189").Value = "y"
$ws.Range("C190").Value = "y"
$ws.Range("D191").Value = "query"
$ws.Range("C192").Value = "y"
$ws.Range("C193").Value = "y"
$ws.Range("C194").Value = "y"
$ws.Range("C195").Value = "y"
$ws.Range("C196").Value = "y"
$ws.Range("C197").Value = "y"
$ws.Range("C198").Value = "y"
$ws.Range("C199").Value = "y"
$ws.Range("C200").Value = "y"
$ws.Range("C201").Value = "y"
$ws.Range("C202").Value = "y"
$ws.Range("C204").Value = "y"
$ws.Range("C205").Value = "y"
$ws.Range("C206").Value = "y"
$ws.Range("C207").Value = "y"
$ws.Range("C208").Value = "y"
$ws.Range("C209").Value = "y"
$ws.Range("C210").Value = "y"
$ws.Range("C211").Value = "ND"
$ws.Range("C212").Value = "y"
$ws.Range("C246").Value = "ND"
$ws.Range("D300").Value = "Query check"
$ws.Range("C305").Value = "y"
$ws.Range("C394").Value = "ND"
$ws.Range("C426").Value = "y"
$ws.Range("C431").Value = "y"
$ws.Range("C437").Value = "y"
$ws.Range("C445").Value = "y"
$ws.Range("C451").Value = "y"
$ws.Range("C458").Value = "y"
$ws.Range("C461").Value = "ND"
$ws.Range("D461").Value = "One of two names of a company"
$ws.Range("C504").Value = "ND"

$ws.Range("C2").Select
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.Zoom = 100
